$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")

# C2 keeps sharing the (updated) timestamp string - its text value is refreshed
# to a newer execution time.
$ws.Range("C2").Value = "Tue Mar 25 16:16:42 IST 2025"

# C3 previously shared the same string as C2; now it gets its own distinct
# (new) timestamp string, one second later.
$ws.Range("C3").Value = "Tue Mar 25 16:16:43 IST 2025"
